# Update gh-pages data output (regenerated crawl numbers)
$wb = $excel.ActiveWorkbook

# ---- Sheet "展览" ----
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 4810
$ws.Range("F3").Value = 2755
$ws.Range("G3").Value = 70
$ws.Range("F5").Value = 2800
$ws.Range("F7").Value = 1954
$ws.Range("F9").Value = 1732
$ws.Range("F10").Value = 753
$ws.Range("F11").Value = 480
$ws.Range("F12").Value = 242
$ws.Range("F13").Value = 411
$ws.Range("F14").Value = 1068
$ws.Range("F18").Value = 544
$ws.Range("F19").Value = 544
$ws.Range("F22").Value = 663
$ws.Range("F23").Value = 743
$ws.Range("F25").Value = 40
$ws.Range("F26").Value = 521
$ws.Range("F27").Value = 26
$ws.Range("F29").Value = 1551
$ws.Range("F30").Value = 347
$ws.Range("F32").Value = 1503
$ws.Range("F33").Value = 190
$ws.Range("F34").Value = 2341
$ws.Range("F35").Value = 393
$ws.Range("F37").Value = 611
$ws.Range("F39").Value = 62
$ws.Range("F41").Value = 788
$ws.Range("F42").Value = 1488
$ws.Range("F46").Value = 43
$ws.Range("F47").Value = 78
$ws.Range("F48").Value = 110

# ---- Sheet "演出" ----
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F2").Value = 8
$ws.Range("F4").Value = 93
$ws.Range("F12").Value = 37

# ---- Sheet "全部类型" ----
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 4810
$ws.Range("F3").Value = 2755
$ws.Range("G3").Value = 70
$ws.Range("F4").Value = 2800
$ws.Range("F5").Value = 1732
$ws.Range("F7").Value = 753
$ws.Range("F8").Value = 480
$ws.Range("F9").Value = 242
$ws.Range("F10").Value = 411
$ws.Range("F11").Value = 1068
$ws.Range("F15").Value = 544
$ws.Range("F16").Value = 544
$ws.Range("F18").Value = 663
$ws.Range("F19").Value = 743
$ws.Range("F21").Value = 93
$ws.Range("F22").Value = 93
$ws.Range("F25").Value = 40
$ws.Range("F26").Value = 521
$ws.Range("F27").Value = 26
$ws.Range("F29").Value = 1551
$ws.Range("F30").Value = 348
$ws.Range("F33").Value = 2341
$ws.Range("F34").Value = 393
$ws.Range("F38").Value = 37
$ws.Range("F39").Value = 612
$ws.Range("F41").Value = 62
$ws.Range("F43").Value = 788
$ws.Range("F44").Value = 1488
$ws.Range("F48").Value = 43
$ws.Range("F49").Value = 78
$ws.Range("F50").Value = 110
